$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) After the title paragraph ("Play Braccio di Ferro Slot Game Free | RTP
#    96.67%"), insert a new "Meta description" paragraph:
#      [empty run] + "Meta description" (bold) + ": Read our review ..." (plain)
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
# Newly split paragraphs inherit the preceding paragraph's style (Heading1);
# reset back to the document default (Normal), matching the source diff,
# which has no explicit pStyle on this paragraph.
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$boldText = "Meta description"
$restText = ": Read our review of Braccio di Ferro slot game and play for free. Enjoy ample opportunities for big wins with a low volatility and high RTP of 96.67%."

$insertionPoint = $d.Range($metaStart, $metaStart)
$insertionPoint.InsertAfter($boldText + $restText)

$boldRange = $d.Range($metaStart, $metaStart + $boldText.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the end of the document, drop the duplicated bold title paragraph
#    ("Play Braccio di Ferro Slot Game Free | RTP 96.67%") entirely, and
#    replace the italic paragraph's text with the new image prompt.
#    Helper: strip the trailing paragraph mark (and any cell-mark char) that
#    Range.Text includes so plain string comparisons work as expected.
# ---------------------------------------------------------------------------
function Get-ParaText($para) {
    return $para.Range.Text.TrimEnd([char]13, [char]7)
}

$dupTitleParaIndex = -1
$promptParaIndex = -1
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $ptext = Get-ParaText $d.Paragraphs($i)
    if ($promptParaIndex -eq -1 -and $ptext -eq "Read our review of Braccio di Ferro slot game and play for free. Enjoy ample opportunities for big wins with a low volatility and high RTP of 96.67%.") {
        $promptParaIndex = $i
    }
    if ($dupTitleParaIndex -eq -1 -and $i -ne 1 -and $ptext -eq "Play Braccio di Ferro Slot Game Free | RTP 96.67%") {
        $dupTitleParaIndex = $i
    }
}

# Delete the duplicated bold title paragraph (whole paragraph incl. mark).
$dupTitlePara = $d.Paragraphs($dupTitleParaIndex)
$dupTitlePara.Range.Delete()

# Re-resolve the prompt paragraph index after the deletion shifted things.
$promptParaIndex2 = -1
$count2 = $d.Paragraphs.Count
for ($i = $count2; $i -ge 1; $i--) {
    $ptext = Get-ParaText $d.Paragraphs($i)
    if ($ptext -eq "Read our review of Braccio di Ferro slot game and play for free. Enjoy ample opportunities for big wins with a low volatility and high RTP of 96.67%.") {
        $promptParaIndex2 = $i
        break
    }
}

$promptPara = $d.Paragraphs($promptParaIndex2)
$pStart = $promptPara.Range.Start
$pEnd = $promptPara.Range.End
$promptRange = $d.Range($pStart, $pEnd)
$promptRange.Text = 'Prompt: Create a feature image for the Braccio di Ferro slot game that captures its adventurous and playful spirit. The image should be in cartoon style and include a happy Maya warrior with glasses. The Maya warrior should be depicted engaging with the game, possibly spinning the reels or collecting flying fish in a bottle. The background should be a small harbor with boats and seagulls in the sky. Use bright colors and playful elements to reflect the fun and excitement of the game. Make sure to include the name of the game, "Braccio di Ferro", in an eye-catching font.'

Write-Output "done"
